$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Fresh" mosquito tape rows (appended after existing data, rows 80-82).
# Values are written column-by-column (D, then E, then A, then B, then C) so
# that the shared-string table gets the new unique strings in the same order
# the original authoring tool produced them.
$ws.Range("D80").Value = "Fresh_F_1"
$ws.Range("D81").Value = "Fresh_F_2"
$ws.Range("D82").Value = "Fresh_F_3"

$ws.Range("E80").Value = "Fresh"
$ws.Range("E81").Value = "Fresh"
$ws.Range("E82").Value = "Fresh"

$ws.Range("A80").Value = "m10"
$ws.Range("A81").Value = "m10"
$ws.Range("A82").Value = "m10"

$ws.Range("B80").Value = 2
$ws.Range("B81").Value = 3
$ws.Range("B82").Value = 4

$ws.Range("C80").Value = 0
$ws.Range("C81").Value = 0
$ws.Range("C82").Value = 0

# Match the centered-alignment style ("s=1") used by every other data row.
$ws.Range("A80:E82").HorizontalAlignment = -4108

# Scroll/selection state as left by the author after the edit.
$ws.Range("A51").Select()
